$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("880×6=5280", $true, $false, $false, $false, $false, $true, 1, $false, "782×8=6256", 2)
$null = $d.Content.Find.Execute("727×9=6543", $true, $false, $false, $false, $false, $true, 1, $false, "360×2=720", 2)
$null = $d.Content.Find.Execute("771×5=3855", $true, $false, $false, $false, $false, $true, 1, $false, "475×4=1900", 2)
$null = $d.Content.Find.Execute("418×8=3344", $true, $false, $false, $false, $false, $true, 1, $false, "720×7=5040", 2)
$null = $d.Content.Find.Execute("415×5=2075", $true, $false, $false, $false, $false, $true, 1, $false, "228×2=456", 2)
$null = $d.Content.Find.Execute("641×5=3205", $true, $false, $false, $false, $false, $true, 1, $false, "403×8=3224", 2)
$null = $d.Content.Find.Execute("581×2=1162", $true, $false, $false, $false, $false, $true, 1, $false, "543×5=2715", 2)
$null = $d.Content.Find.Execute("106×3=318", $true, $false, $false, $false, $false, $true, 1, $false, "946×3=2838", 2)
$null = $d.Content.Find.Execute("521×4=2084", $true, $false, $false, $false, $false, $true, 1, $false, "991×7=6937", 2)
$null = $d.Content.Find.Execute("835×9=7515", $true, $false, $false, $false, $false, $true, 1, $false, "204×6=1224", 2)
$null = $d.Content.Find.Execute("690×8=5520", $true, $false, $false, $false, $false, $true, 1, $false, "799×8=6392", 2)
$null = $d.Content.Find.Execute("147×2=294", $true, $false, $false, $false, $false, $true, 1, $false, "370×8=2960", 2)
$null = $d.Content.Find.Execute("535×4=2140", $true, $false, $false, $false, $false, $true, 1, $false, "453×3=1359", 2)
$null = $d.Content.Find.Execute("194×2=388", $true, $false, $false, $false, $false, $true, 1, $false, "982×9=8838", 2)
$null = $d.Content.Find.Execute("405×5=2025", $true, $false, $false, $false, $false, $true, 1, $false, "857×7=5999", 2)
$null = $d.Content.Find.Execute("597×7=4179", $true, $false, $false, $false, $false, $true, 1, $false, "438×4=1752", 2)
$null = $d.Content.Find.Execute("920×7=6440", $true, $false, $false, $false, $false, $true, 1, $false, "568×7=3976", 2)
$null = $d.Content.Find.Execute("770×5=3850", $true, $false, $false, $false, $false, $true, 1, $false, "726×8=5808", 2)
$null = $d.Content.Find.Execute("656×9=5904", $true, $false, $false, $false, $false, $true, 1, $false, "911×5=4555", 2)
$null = $d.Content.Find.Execute("314×6=1884", $true, $false, $false, $false, $false, $true, 1, $false, "735×6=4410", 2)
$null = $d.Content.Find.Execute("804×7=5628", $true, $false, $false, $false, $false, $true, 1, $false, "617×8=4936", 2)
$null = $d.Content.Find.Execute("309×5=1545", $true, $false, $false, $false, $false, $true, 1, $false, "236×9=2124", 2)
$null = $d.Content.Find.Execute("355×8=2840", $true, $false, $false, $false, $false, $true, 1, $false, "169×7=1183", 2)
$null = $d.Content.Find.Execute("697×9=6273", $true, $false, $false, $false, $false, $true, 1, $false, "162×2=324", 2)
$null = $d.Content.Find.Execute("216×9=1944", $true, $false, $false, $false, $false, $true, 1, $false, "734×8=5872", 2)
